$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must remain plain text (some contain
# multiple "." separators or leading/trailing zeros that a numeric cast would
# mangle), so force a Text number format on each D cell before writing it.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.090.88"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.047.82"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.76"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.91"
$ws.Range("E8").Value = "  -4.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.87"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.865"
$ws.Range("E13").Value = "  +5.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.346.56"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.65"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.047.83"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.83"
$ws.Range("E17").Value = "  +14.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.028.01"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.60"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -5.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.35"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.63"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.42"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.74"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.96"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.79"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0615"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0887"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.23"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.33"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.17"
$ws.Range("E39").Value = "  +12.56%  "
$ws.Range("E40").Value = "  +15.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0221"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.11"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.13"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.48"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0927"
$ws.Range("E45").Value = "  -21.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.41"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.268.06"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.79"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.229.79"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.50"
$ws.Range("E51").Value = "  -1.36%  "
